# New advice filled in (Sheet "A" -> person "Romviel").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Last name of the person on sheet 1 was a placeholder ("A"); fill in the real name.
$ws.Range("B2").Value = "Romviel"

# The "Voor het ontbijt" advice cell had a leftover placeholder value; clear it.
$ws.Range("C11").Value = ""

# New nutrition ("Voedingsadvies") advice text for each moment of the day.
$ws.Range("C14").Value = "2 volkoren boterhammen, 2 porties eiwitrijk beleg, 200 ml sojamelk`n"
$ws.Range("C15").Value = "1,5 maatschepje PROMISS eiwitpoeder, 30 gram ongezouten noten, 50 gram rauwkost"
$ws.Range("C17").Value = "100 gram fruit of 2 blokjes pure chocolade "
$ws.Range("C16").Value = "200 gram groente, 70 gram onbereid vlees of 50 gram onbereid vis of 100 gram onbereid vegetarisch product, 125 gram aardappelen of 75 gram rijst of 50 gram pasta "
$ws.Range("C13").Value = "1,5 maatschepje PROMISS eiwitpoeder, 30 gram ongezouten noten"
$ws.Range("C12").Value = "Griekse yoghurt 10%, 4 dadels, 4 lepels (50 gram) muesli, 200 ml sojamelk"

# The lunch advice is now longer, so wrap the text in that cell.
$ws.Range("C14").WrapText = $true

# Update the view: zoomed in a bit more and selection moved to B18.
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("B18").Select()
